$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New loading_percent values for rows 2-25 (data rows 0-23), columns B,D,E,F,G,H,I,K,L,O
$data = @(
    @(13.21241398033335, 4.490498767714787, 17.20030952492571, 21.12345547573798, 23.28252022839287, 12.96501463178743, 24.77110034786419, 11.12605332938826, 8.54489681819906, 18.99301785127283),
    @(13.03818763570395, 4.402932019123607, 17.2787578360615, 21.15935016798435, 23.34140368521217, 13.0106743175468, 24.90736185898097, 10.82624240069859, 8.492276400498444, 19.06587109879511),
    @(12.93194405226256, 4.347602460565738, 17.32980918305589, 21.18750770153787, 23.38632357870721, 13.04082007736894, 24.99559138265353, 10.63651236579948, 8.461144381087818, 19.11495887680282),
    @(12.88888071765423, 4.324678765969831, 17.35133944057484, 21.20051723370785, 23.40682269565091, 13.05363550102815, 25.03269560573958, 10.55785340355105, 8.448764038926722, 19.13605560710621),
    @(12.88174547321763, 4.320850069756565, 17.3549584350551, 21.20277006281009, 23.41035874405378, 13.05579555252099, 25.03892625517433, 10.54471342829479, 8.446727076336575, 19.13962465634951),
    @(12.93136228369578, 4.347294805253753, 17.33009660468355, 21.18767694194354, 23.38659116828575, 13.04099076144015, 24.99608712400449, 10.63545687294561, 8.460976162479717, 19.11523897185208),
    @(13.15221743882296, 4.460639913913456, 17.22676074226898, 21.13456036988794, 23.30099824796785, 12.9803200955735, 24.81713731100841, 11.02389452619027, 8.526515241346244, 19.01723235288344),
    @(13.58902803684631, 4.669819965841325, 17.04694294169067, 21.07906179571648, 23.2030993678748, 12.87808782791969, 24.50233189607874, 11.73774373284824, 8.663909466286038, 18.85969675591524),
    @(13.90943534841954, 4.814670939055249, 16.92866391331452, 21.06807121041151, 23.17427944351096, 12.81317976189814, 24.29291742587488, 12.22930079792976, 8.769610737863239, 18.76520899280798),
    @(14.05450837723883, 4.878486696113764, 16.87784259919292, 21.06954786628186, 23.17059934328208, 12.78586578329877, 24.20236891631494, 12.44510275707395, 8.818585105468379, 18.72686491540225),
    @(14.10929899853154, 4.902340876604271, 16.85902583213932, 21.07103768073378, 23.17056528895258, 12.77584090853747, 24.16875646227642, 12.52564905374923, 8.837245884843476, 18.7130141084932),
    @(14.09750606242204, 4.897217496849652, 16.86305933254049, 21.07067545161456, 23.17051211934396, 12.77798578487747, 24.17596545409918, 12.50835489262538, 8.833222031476105, 18.71596732793238),
    @(14.05901921443359, 4.880455505541535, 16.87628595704212, 21.06965178961019, 23.17056928042464, 12.78503465026956, 24.1995900551753, 12.45175313205709, 8.820118102996357, 18.72571197589727),
    @(14.03542459459236, 4.870147379116806, 16.88444337671031, 21.06914592966765, 23.17078141948235, 12.78939374381824, 24.21414882919177, 12.41692872793298, 8.812106195351339, 18.73176807558739),
    @(13.89993668581794, 4.810457515152637, 16.93204517421922, 21.06810499143472, 23.17471003807558, 12.81500933532485, 24.29892971593084, 12.21503607366964, 8.766426974840877, 18.76780840454052),
    @(13.81661048669787, 4.773298252654123, 16.96201099962094, 21.06912495774588, 23.1795383749388, 12.83129055405513, 24.35214646458568, 12.08914386682014, 8.738623249406386, 18.79110770018502),
    @(13.76862203073703, 4.751730367566717, 16.9795275143423, 21.07032129043094, 23.18320301284857, 12.84086337435262, 24.38319924698665, 12.01600157024621, 8.722715860870521, 18.80494534632952),
    @(13.75236487833944, 4.74439478656943, 16.98550659183445, 21.07083106638559, 23.18459609577817, 12.84414033835218, 24.39378948227862, 11.99111266229765, 8.717344805566968, 18.80970543768412),
    @(13.8254874079589, 4.777274181705728, 16.95879201485484, 21.06895328353783, 23.17893250671318, 12.82953583213296, 24.34643551992101, 12.10262149031282, 8.741574341892273, 18.78858225727988),
    @(14.07032805051274, 4.885387458352667, 16.87238936441236, 21.06992721613431, 23.17051557385649, 12.78295558722417, 24.19263259613923, 12.46841065708303, 8.82396401921025, 18.72283155591743),
    @(14.22947723190309, 4.954225045220657, 16.81841544845947, 21.07598729824263, 23.17293945924281, 12.75436826276427, 24.09605472573297, 12.70061679701032, 8.878476171965717, 18.68376160129391),
    @(14.14463103126432, 4.91765568289436, 16.84699435270325, 21.07225709709345, 23.17091991660043, 12.7694560362568, 24.14724009070945, 12.57732683105911, 8.849325354238873, 18.70425624284621),
    @(13.82147440712022, 4.775477300277095, 16.9602464183342, 21.06902899750129, 23.17920365160408, 12.83032847990064, 24.34901601196747, 12.0965306362919, 8.740239911096936, 18.78972263160792),
    @(13.47075109863154, 4.614724217584811, 17.09315344093987, 21.08885068997491, 23.22204193512896, 12.90395246255083, 24.5836437745693, 11.55016042413747, 8.625856735479392, 18.89859085035387)
)

$colIndices = @(2, 4, 5, 6, 7, 8, 9, 11, 12, 15)  # B, D, E, F, G, H, I, K, L, O

$startRow = 2
for ($r = 0; $r -lt $data.Count; $r++) {
    $rowData = $data[$r]
    for ($c = 0; $c -lt $colIndices.Count; $c++) {
        $ws.Cells.Item($startRow + $r, $colIndices[$c]).Value = $rowData[$c]
    }
}
